# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across the resume.
#
# Strategy: for each paragraph that contains metrics, walk left-to-right
# through the paragraph's Range, Find()-ing each metric substring in turn
# (search scope collapses forward after each hit, just like real Word COM),
# and set Font.Bold + Font.Color on the found sub-range. This naturally
# splits the paragraph's single run into multiple runs exactly like the
# target OOXML diff.

$d = $word.ActiveDocument
$highlightColor = 5258796   # BGR-packed wdColor for RGB(0x2C,0x3E,0x50) == "2C3E50"

function Highlight-Metrics {
    param(
        [int]$ParaIndex,
        [string[]]$Targets
    )

    $para = $d.Paragraphs($ParaIndex)
    $cursor = $para.Range.Start
    $paraEnd = $para.Range.End

    foreach ($target in $Targets) {
        $searchRange = $d.Range($cursor, $paraEnd)
        $found = $searchRange.Find.Execute($target, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if ($found) {
            $searchRange.Font.Bold = $true
            $searchRange.Font.Color = $highlightColor
            $cursor = $searchRange.End
        }
    }
}

# • Discovered systematic race coding errors ... from 23% to 64%
Highlight-Metrics 10 @("23%", "64%")

# • Utilized advanced sampling methods ... from ±4.2% to ±2.1%, ... from 71% to 87%, ...
Highlight-Metrics 12 @("±4.2%", "±2.1%", "71%", "87%")

# • Trigonometric algorithm ... reduced mapping costs by 73.5%, saving ... $4.7M ...
Highlight-Metrics 13 @("73.5%", "`$4.7M")

# • Built real-time FEC analysis systems ... valued over $2 trillion
Highlight-Metrics 14 @("`$2")

# • Modernized legacy ETL processes ... reducing processing time by 57%
Highlight-Metrics 20 @("57%")

# • 178% accuracy improvement in racial classification algorithms
Highlight-Metrics 85 @("178%")

# • Algorithmic innovation: ... reducing mapping costs 73.5%
Highlight-Metrics 86 @("73.5%")

# • $4.7M savings enabled nonprofit access
Highlight-Metrics 87 @("`$4.7M")

# • Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations
Highlight-Metrics 88 @("12,847")

# • Predictive excellence: ... margin of error from ±4.2% to ±2.1%
Highlight-Metrics 90 @("±4.2%", "±2.1%")

# • Increased voter turnout prediction accuracy from 71% to 87%
Highlight-Metrics 91 @("71%", "87%")

Write-Host "Highlighting complete"
